$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "GRT-USD"
$ws.Range("A25").Value = "BSCX-USD"
